# "Generate Report for Handback"
# Updates the in-flight handoff (3ae5b4a1...) to its new guid (a9dd0855...)
# with refreshed timestamps, and appends a brand-new handback row
# (fcde77c9...) to the Overview, zh-cn and de-de report sheets.

$wb = $excel.ActiveWorkbook

$oldGuid = "3ae5b4a1-f782-4e07-859a-f90961d9d7f3"
$newGuid = "a9dd0855-3691-45aa-ad81-64807975381c"
$addGuid = "fcde77c9-9765-4063-aae6-a81d4ec2b71a"

$zhOldHash = "d05692d9e3904b7f8bc0ed4efe5f27ca373d7998"
$zhNewHash = "a3080b6a1099c53ddbb5b815125a03d8c29fb6ab"
$deOldHash = "d05692d9e3904b7f8bc0ed4efe5f27ca373d7998"
$deNewHash = "a3080b6a1099c53ddbb5b815125a03d8c29fb6ab"
$addHash   = "21f4d94d58f71496b203d73608e95bfdfa119546"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "$newGuid.md"
$ws1.Range("B2").Value = "e2e\$newGuid.md"
$ws1.Range("G2").Value = "2016-09-06 01:08:14"

$lo1 = $ws1.ListObjects.Item(1)
$lo1.ListRows.Add() | Out-Null

$ws1.Range("A3").Value = "$addGuid.md"
$ws1.Range("B3").Value = "e2e\$addGuid.md"
$ws1.Range("C3").Value = ".md"
$ws1.Range("E3").Value = "Handed back: in sync with en-US"
$ws1.Range("F3").Value = "Handed back: in sync with en-US"
$ws1.Range("G3").Value = "2016-09-06 01:08:14"

$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6d76bb861347f8643d0c2f57e375745822ac91fe/e2e/$addGuid.md", "", "", "e2e\$addGuid.md") | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "$newGuid.md"
$ws2.Range("G2").Value = "$newGuid.$zhNewHash.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-09-06 01:08:03"
$ws2.Range("I2").Value = "$newGuid.md"
$ws2.Range("J2").Value = "$newGuid.$zhNewHash.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-09-06 01:08:28"

$lo2 = $ws2.ListObjects.Item(1)
$lo2.ListRows.Add() | Out-Null

$ws2.Range("A3").Value = "$addGuid.md"
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Handed back: in sync with en-US"
$ws2.Range("D3").Value = "e2e"
$ws2.Range("E3").Value = "ht"
$ws2.Range("F3").Value = "True"
$ws2.Range("G3").Value = "$addGuid.$addHash.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-09-06 01:08:03"
$ws2.Range("I3").Value = "$addGuid.md"
$ws2.Range("J3").Value = "$addGuid.$addHash.zh-cn.xlf"
$ws2.Range("K3").Value = "2016-09-06 01:08:28"
$ws2.Range("L3").Value = ""
$ws2.Range("M3").Value = "True"
$ws2.Range("N3").Value = ""
$ws2.Range("O3").Value = "False"
$ws2.Range("P3").Value = ""

$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6d76bb861347f8643d0c2f57e375745822ac91fe/e2e/$addGuid.md", "", "", "$addGuid.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/78fede49d63bcbfb2dc66b1deab1b509bbdc5c81/e2e/$addGuid.md", "", "", "$addGuid.md") | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "$newGuid.md"
$ws3.Range("G2").Value = "$newGuid.$deNewHash.de-de.xlf"
$ws3.Range("H2").Value = "2016-09-06 01:08:14"
$ws3.Range("I2").Value = "$newGuid.md"
$ws3.Range("J2").Value = "$newGuid.$deNewHash.de-de.xlf"
$ws3.Range("K2").Value = "2016-09-06 01:08:37"

$lo3 = $ws3.ListObjects.Item(1)
$lo3.ListRows.Add() | Out-Null

$ws3.Range("A3").Value = "$addGuid.md"
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Handed back: in sync with en-US"
$ws3.Range("D3").Value = "e2e"
$ws3.Range("E3").Value = "ht"
$ws3.Range("F3").Value = "True"
$ws3.Range("G3").Value = "$addGuid.$addHash.de-de.xlf"
$ws3.Range("H3").Value = "2016-09-06 01:08:14"
$ws3.Range("I3").Value = "$addGuid.md"
$ws3.Range("J3").Value = "$addGuid.$addHash.de-de.xlf"
$ws3.Range("K3").Value = "2016-09-06 01:08:37"
$ws3.Range("L3").Value = ""
$ws3.Range("M3").Value = "True"
$ws3.Range("N3").Value = ""
$ws3.Range("O3").Value = "False"
$ws3.Range("P3").Value = ""

$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6d76bb861347f8643d0c2f57e375745822ac91fe/e2e/$addGuid.md", "", "", "$addGuid.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/2442417c3040476793cc81a92e433fe1d30f6396/e2e/$addGuid.md", "", "", "$addGuid.md") | Out-Null

Write-Host "Handback report regenerated"
